$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: shift the two existing data rows down by one (row2->3,
#    row3->4), then insert the remaining rows needed for a 33-row table.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item("5:33").Insert()

# ---------------------------------------------------------------------
# 2. Full set of question rows (row -> [A, B, C]). Row 2 is brand new and
#    keeps the plain default column styles (A=10, B=1); only column C
#    needs the bordered look applied further down.
# ---------------------------------------------------------------------
$data = @{}
$data[2] = @("Введите название вашего стартапа", "name_project", "text")
$data[3] = @("Описание краткой истории компании и ее деятельности (1-2 предложения)", "product_name", "text")
$data[4] = @("Почему в вас стоит инвестировать?`nПродукт (фокусный продукт, на развитие которого привлекаются инвестиции) (1-2 предложения)", "feature", "text")
$data[5] = @("Краткое описание Вашего стартапа (2 - 3 предложения)", "about_project", "text")
$data[6] = @("Опишите основной эффект от использования Вашей идеи", "еffectiveness", "text")
$data[7] = @("Выделите 2-3 основные цели Вашего стартапа", "goals", "text")
$data[8] = @("Возможные негативные эффекты от Вашего стартапа", "risks", "text")
$data[9] = @("Сколько примерно человек в штате? (в следующих вопросах Вам будет необходимо их представить)", "team", "text")
$data[10] = @("Фамилия И.О., Должность, Опыт", "teammate1", "text")
$data[11] = @("Фотография сотрудника", "teammate1_photo", "image")
$data[12] = @("Фамилия И.О., Должность, Опыт  (чтобы пропустить вопрос напишите `"-`")", "teammate2", "text")
$data[13] = @("Фотография сотрудника (чтобы пропустить вопрос напишите `"-`")", "teammate2_photo", "image")
$data[14] = @("Фамилия И.О., Должность, Опыт  (чтобы пропустить вопрос напишите `"-`")", "teammate3", "text")
$data[15] = @("Фотография сотрудника (чтобы пропустить вопрос напишите `"-`")", "teammate3_photo", "image")
$data[16] = @("Фамилия И.О., Должность, Опыт  (чтобы пропустить вопрос напишите `"-`")", "teammate4", "text")
$data[17] = @("Фотография сотрудника (чтобы пропустить вопрос напишите `"-`")", "teammate4_photo", "image")
$data[18] = @("Напишите примерную дату первого этапа Вашего стартапа (Всего этапов будет 4)", "first_stage", "text")
$data[19] = @("Опишите (пару слов) о данном этапе", "first_stage_description", "text")
$data[20] = @("Напишите примерную дату второго этапа Вашего стартапа", "second_stage", "text")
$data[21] = @("Опишите (пару слов) о данном этапе", "second_stage_discription", "text")
$data[22] = @("Напишите примерную дату третьего этапа Вашего стартапа", "third_stage", "text")
$data[23] = @("Опишите (пару слов) о данном этапе", "third_stage_discription", "text")
$data[24] = @("Напишите примерную дату заключающего этапа Вашего стартапа", "final_stage", "text")
$data[25] = @("Опишите конечный результат Вашего стартапа", "final_stage_discription", "text")
$data[26] = @("Напишите актуальность Вашего сартапа", "problems", "text")
$data[27] = @("Инструменты для реализации Вашего решения (пару слов)", "decision", "text")
$data[28] = @("Что даёт возможность реализовать Ваш стартап?", "strategy", "text")
$data[29] = @("Каков прогноз Вашего стартапа? (фотография где сравниваются Ваша стратегия с текущим положением рынка)", "forecast", "image")
$data[30] = @("Загрузите первую фотографию Вашего решения (Всего фотографий 2)", "picture1", "image")
$data[31] = @("Загрузите вторую фотографию Вашего решения", "picture2", "image")
$data[32] = @("Напишите описание к второй фотографии", "descript_picture1", "text")
$data[33] = @("Напишите описание к первой фотографии", "descript_picture2", "text")


foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------
# 3. Column C uses the "text/image" bordered style (as column C2 already
#    does) for every data row - copy that format down the whole column.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C2:C33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. A handful of column-A cells use the "section header" look instead
#    of the plain style - row 4 (original row3 template) already has it
#    after the row-insert shift; copy it onto row 9 as well.
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial can carry the source text/style, so
# make sure the final text is still correct after copying formats).
$ws.Cells.Item(9, 1).Value = $data[9][0]

# ---------------------------------------------------------------------
# 5. Row heights that Excel widened for wrapped text.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(22).RowHeight = 29.4
$ws.Rows.Item(23).RowHeight = 29.4

# ---------------------------------------------------------------------
# 6. AutoFilter over the full table + the matching hidden defined name
#    Excel writes alongside it.
# ---------------------------------------------------------------------
$ws.Range("A1:C33").AutoFilter()
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Лист1!`$A`$1:`$C`$33")
$nm.Visible = $false

# ---------------------------------------------------------------------
# 7. View state: zoom + selection + window placement/size.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("G12").Select()
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 17256
$excel.ActiveWindow.Height = 5928

Write-Host "edit complete"
